# Update frequency table values for publication run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00608519269776876
$ws.Range("C2").Value = 0.979716024340771
$ws.Range("D2").Value = 0.00405679513184584
$ws.Range("E2").Value = 0.00202839756592292
$ws.Range("F2").Value = 0.0162271805273834
$ws.Range("G2").Value = 0.0101419878296146
$ws.Range("H2").Value = 0.0121703853955375
$ws.Range("I2").Value = 0.971602434077079
$ws.Range("J2").Value = 0.0202839756592292
$ws.Range("K2").Value = 0.947261663286004
$ws.Range("L2").Value = 0.961460446247465
$ws.Range("M2").Value = 0.00202839756592292
$ws.Range("N2").Value = 0.00608519269776876
$ws.Range("Q2").Value = 0.995943204868154
$ws.Range("R2").Value = 0.00405679513184584
$ws.Range("S2").Value = 0.886409736308316
$ws.Range("T2").Value = 0.0567951318458418
$ws.Range("U2").Value = 0.0689655172413793
$ws.Range("V2").Value = 0.987829614604463
$ws.Range("W2").Value = 0.00811359026369168
$ws.Range("X2").Value = 0.00202839756592292
$ws.Range("B3").Value = 0.0121703853955375
$ws.Range("C3").Value = 0.00202839756592292
$ws.Range("D3").Value = 0.0223123732251521
$ws.Range("E3").Value = 0.00608519269776876
$ws.Range("F3").Value = 0.00405679513184584
$ws.Range("H3").Value = 0.975659229208925
$ws.Range("I3").Value = 0.00608519269776876
$ws.Range("J3").Value = 0.0162271805273834
$ws.Range("K3").Value = 0.00405679513184584
$ws.Range("L3").Value = 0.00202839756592292
$ws.Range("M3").Value = 0.941176470588235
$ws.Range("N3").Value = 0.00811359026369168
$ws.Range("P3").Value = 0.993914807302231
$ws.Range("Q3").Value = 0.00405679513184584
$ws.Range("R3").Value = 0.993914807302231
$ws.Range("S3").Value = 0.0993914807302231
$ws.Range("T3").Value = 0.935091277890467
$ws.Range("U3").Value = 0.920892494929006
$ws.Range("W3").Value = 0.00202839756592292
$ws.Range("X3").Value = 0.00405679513184584
$ws.Range("B4").Value = 0.979716024340771
$ws.Range("C4").Value = 0.00405679513184584
$ws.Range("D4").Value = 0.0121703853955375
$ws.Range("E4").Value = 0.00608519269776876
$ws.Range("F4").Value = 0.977687626774848
$ws.Range("G4").Value = 0.987829614604463
$ws.Range("H4").Value = 0.00811359026369168
$ws.Range("I4").Value = 0.00608519269776876
$ws.Range("J4").Value = 0.955375253549696
$ws.Range("K4").Value = 0.00405679513184584
$ws.Range("L4").Value = 0.0344827586206897
$ws.Range("M4").Value = 0.00202839756592292
$ws.Range("N4").Value = 0.00202839756592292
$ws.Range("P4").Value = 0.00202839756592292
$ws.Range("R4").Value = 0.00202839756592292
$ws.Range("S4").Value = 0.0101419878296146
$ws.Range("T4").Value = 0.00202839756592292
$ws.Range("V4").Value = 0.0101419878296146
$ws.Range("W4").Value = 0.98580121703854
$ws.Range("X4").Value = 0.987829614604463
$ws.Range("B5").Value = 0.00202839756592292
$ws.Range("C5").Value = 0.0141987829614604
$ws.Range("D5").Value = 0.961460446247465
$ws.Range("E5").Value = 0.98580121703854
$ws.Range("F5").Value = 0.00202839756592292
$ws.Range("G5").Value = 0.00202839756592292
$ws.Range("H5").Value = 0.00405679513184584
$ws.Range("I5").Value = 0.0162271805273834
$ws.Range("J5").Value = 0.00811359026369168
$ws.Range("K5").Value = 0.0425963488843813
$ws.Range("M5").Value = 0.0527383367139959
$ws.Range("N5").Value = 0.983772819472617
$ws.Range("P5").Value = 0.00405679513184584
$ws.Range("S5").Value = 0.00405679513184584
$ws.Range("T5").Value = 0.00405679513184584
$ws.Range("U5").Value = 0.0101419878296146
$ws.Range("V5").Value = 0.00202839756592292
$ws.Range("W5").Value = 0.00405679513184584
$ws.Range("X5").Value = 0.00608519269776876
